# Update "想去人数" (want-to-go count, column F) figures in the
# "展览" (Exhibition) and "全部类型" (All types) sheets to match the
# newly generated gh-pages data snapshot.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

# --- Sheet "展览" (sheet1) ---
$ws1.Range("F3").Value = 28
$ws1.Range("F5").Value = 15966
$ws1.Range("F6").Value = 422
$ws1.Range("F8").Value = 719
$ws1.Range("F9").Value = 15504
$ws1.Range("F11").Value = 9129
$ws1.Range("F12").Value = 403
$ws1.Range("F16").Value = 209
$ws1.Range("F18").Value = 211
$ws1.Range("F20").Value = 70
$ws1.Range("F24").Value = 68
$ws1.Range("F25").Value = 1127
$ws1.Range("F26").Value = 9
$ws1.Range("F27").Value = 18
$ws1.Range("F28").Value = 26
$ws1.Range("F29").Value = 501
$ws1.Range("F33").Value = 69
$ws1.Range("F34").Value = 55
$ws1.Range("F35").Value = 264
$ws1.Range("F36").Value = 334
$ws1.Range("F38").Value = 119
$ws1.Range("F39").Value = 5614
$ws1.Range("F40").Value = 5235

# --- Sheet "全部类型" (sheet4) ---
$ws4.Range("F3").Value = 28
$ws4.Range("F5").Value = 15967
$ws4.Range("F6").Value = 422
$ws4.Range("F8").Value = 719
$ws4.Range("F9").Value = 15504
$ws4.Range("F11").Value = 9129
$ws4.Range("F12").Value = 403
$ws4.Range("F16").Value = 209
$ws4.Range("F18").Value = 211
$ws4.Range("F20").Value = 70
$ws4.Range("F24").Value = 68
$ws4.Range("F25").Value = 1127
$ws4.Range("F26").Value = 9
$ws4.Range("F27").Value = 18
$ws4.Range("F28").Value = 26
$ws4.Range("F29").Value = 501
$ws4.Range("F35").Value = 69
$ws4.Range("F36").Value = 55
$ws4.Range("F37").Value = 264
$ws4.Range("F38").Value = 334
$ws4.Range("F40").Value = 119
$ws4.Range("F41").Value = 5614
$ws4.Range("F43").Value = 5235
